# Update the NATMI Bmp4-Bmpr1b LR-pairs sheet with refreshed TPM-derived
# values, and drop the "Resolving-Mac" target cluster rows entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old data rows (2-10) completely so stale cells/rows beyond the
# new data range (which only needs rows 2-7) are removed.
$ws.Range("A2:T10").EntireRow.Delete() | Out-Null

# New data, one row per sending/target cluster pair (Resolving-Mac removed).
$data = @(
    @("ECs",  "Bmp4", "Bmpr1b", "FAPs",  3, 1, 6.325956666666666, 18.97787,  0.4468357575736242, 0.4592138460625664, 3, 1, 1.467406, 4.402218, 0.864087546066766,  0.9050919696083439, 9.282746768406664, 83.54472091565999, 0.3861052132566773, 0.4156307644041911),
    @("ECs",  "Bmp4", "Bmpr1b", "MuSCs", 3, 1, 6.325956666666666, 18.97787,  0.4468357575736242, 0.4592138460625664, 2, 1, 0.2308085, 0.461617, 0.135912453933234, 0.09490803039165596, 1.460084569298333, 8.76050741579,      0.06073054431694692, 0.04358308165837527),
    @("FAPs", "Bmp4", "Bmpr1b", "FAPs",  3, 1, 6.686451000000001, 20.059353, 0.4722993778644153, 0.4853828506917099, 3, 1, 1.467406, 4.402218, 0.864087546066766,  0.9050919696083439, 9.811738316106,    88.305644844954,    0.4081080104277229, 0.4393161203466724),
    @("FAPs", "Bmp4", "Bmpr1b", "MuSCs", 3, 1, 6.686451000000001, 20.059353, 0.4722993778644153, 0.4853828506917099, 2, 1, 0.2308085, 0.461617, 0.135912453933234, 0.09490803039165596, 1.5432897256335,   9.259738353801001,  0.06419136743669243, 0.04606673034503741),
    @("MuSCs","Bmp4", "Bmpr1b", "FAPs",  2, 1, 1.1448225,          2.289645,  0.08086486456196039, 0.05540330324572383, 3, 1, 1.467406, 4.402218, 0.864087546066766,  0.9050919696083439, 1.679919405435,    10.07951643261,     0.06987432238236574, 0.05014508485748054),
    @("MuSCs","Bmp4", "Bmpr1b", "MuSCs", 2, 1, 1.1448225,          2.289645,  0.08086486456196039, 0.05540330324572383, 2, 1, 0.2308085, 0.461617, 0.135912453933234, 0.09490803039165596, 0.26423476399125,  1.056939055965,     0.01099054217959465, 0.005258218388243288)
)

$r = 2
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
